# Apply the crypto-list refresh described by the commit:
# "Updated cryptos list on Fri Mar 10 18:41:08 UTC 2023 with GitHub Actions"
#
# For every affected row we overwrite the Price (D) and/or Volume(1h) (E)
# text, and for two rows (44/45) the Coin name (B) and Link (C) as well -
# those two rows swapped which coin occupies which row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a plain text value (Coin / Link / Volume columns). These
# values are never numeric-looking, so a normal .Value assignment already
# keeps them as text.
function Set-TextValue($addr, $val) {
    $ws.Range($addr).Value = $val
}

# Helper: write a text value into the Price column (D). Several of these
# look like plain numbers ("1.003", "275.55", ...) which Excel would
# otherwise auto-convert to a numeric value. Force the cell to Text format
# first, then restore the default "Normal" style so no stray formatting is
# left behind on the cell.
function Set-PriceValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# row -> updated values (only columns that actually changed are listed)
$rowUpdates = @(
    @{ Row = 2;  D = "19.935.55";  E = "  -5.92%  " },
    @{ Row = 3;  D = "1.408.68";   E = "  -7.29%  " },
    @{ Row = 4;  D = "1.003";      E = "  -0.76%  " },
    @{ Row = 5;                    E = "  -0.25%  " },
    @{ Row = 6;  D = "275.55";     E = "  -4.20%  " },
    @{ Row = 7;  D = "0.3655";     E = "  -6.00%  " },
    @{ Row = 8;  D = "0.3105";     E = "  -1.59%  " },
    @{ Row = 9;  D = "39.83";      E = "  -6.55%  " },
    @{ Row = 10;                   E = "  -2.79%  " },
    @{ Row = 11; D = "0.06498";    E = "  -8.06%  " },
    @{ Row = 12;                   E = "  -0.88%  " },
    @{ Row = 13; D = "5.502";      E = "  -2.94%  " },
    @{ Row = 14;                   E = "  -1.97%  " },
    @{ Row = 15; D = "6.179";      E = "  -3.93%  " },
    @{ Row = 16; D = "1.410.93";   E = "  -7.55%  " },
    @{ Row = 17; D = "0.00001018"; E = "  -5.63%  " },
    @{ Row = 18; D = "0.05674";    E = "  -14.19%  " },
    @{ Row = 19; D = "1.002";      E = "  -0.07%  " },
    @{ Row = 20; D = "70.79";      E = "  -14.30%  " },
    @{ Row = 21; D = "5.610";      E = "  -7.44%  " },
    @{ Row = 22; D = "14.71";      E = "  -3.57%  " },
    @{ Row = 23; D = "10.90";      E = "  +1.50%  " },
    @{ Row = 24; D = "2.241";      E = "  -5.69%  " },
    @{ Row = 25; D = "19.958.88";  E = "  -5.87%  " },
    @{ Row = 26; D = "2.257";      E = "  -4.70%  " },
    @{ Row = 27; D = "132.74";     E = "  -10.31%  " },
    @{ Row = 28; D = "17.30";      E = "  -5.27%  " },
    @{ Row = 29; D = "1.569.46";   E = "  -7.57%  " },
    @{ Row = 30; D = "109.76";     E = "  -5.51%  " },
    @{ Row = 31; D = "3.967";      E = "  -17.34%  " },
    @{ Row = 32; D = "5.275";      E = "  -12.16%  " },
    @{ Row = 33; D = "0.8150";     E = "  -14.34%  " },
    @{ Row = 34; D = "0.07689";    E = "  -4.02%  " },
    @{ Row = 35; D = "8.318";      E = "  -1.38%  " },
    @{ Row = 36; D = "1.475";      E = "  -1.50%  " },
    @{ Row = 37; D = "4.910";      E = "  -3.92%  " },
    @{ Row = 38; D = "0.05855";    E = "  -0.96%  " },
    @{ Row = 39;                   E = "  -0.31%  " },
    @{ Row = 40; D = "0.02061";    E = "  -4.55%  " },
    @{ Row = 41; D = "10.48";      E = "  -6.91%  " },
    @{ Row = 42; D = "0.1897";     E = "  -5.13%  " },
    @{ Row = 43;                   E = "  -6.17%  " },
    @{ Row = 44; B = "TheSandbox"; C = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"; D = "0.5305"; E = "  -6.76%  " },
    @{ Row = 45; B = "EnergySwap"; C = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens";   D = "12.35";  E = "  -5.04%  " },
    @{ Row = 46;                   E = "  -4.64%  " },
    @{ Row = 47; D = "0.5182";     E = "  -5.91%  " },
    @{ Row = 48; D = "114.74";     E = "  -0.31%  " },
    @{ Row = 49; D = "1.769";      E = "  -5.73%  " },
    @{ Row = 50;                   E = "  -9.87%  " },
    @{ Row = 51;                   E = "  -0.83%  " }
)

foreach ($update in $rowUpdates) {
    $row = $update.Row

    if ($update.ContainsKey("B")) {
        Set-TextValue "B$row" $update.B
    }
    if ($update.ContainsKey("C")) {
        Set-TextValue "C$row" $update.C
    }
    if ($update.ContainsKey("D")) {
        Set-PriceValue "D$row" $update.D
    }
    if ($update.ContainsKey("E")) {
        Set-TextValue "E$row" $update.E
    }
}
